$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24
$srcRow = 23

# Clone the formatting of the previous data row onto the new row so the
# styled columns (A: bold/border/center index, E: date-time format) match.
$ws.Range("A" + $srcRow + ":V" + $srcRow).Copy()
$ws.Range("A" + $row + ":V" + $row).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 23
$ws.Cells.Item($row, 2).Value = "india"
$ws.Cells.Item($row, 3).Value = "isl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45225.6875
$ws.Cells.Item($row, 6).Value = "North East Utd"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Jamshedpur"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 3.24
$ws.Cells.Item($row, 11).Value = "22/10/2023 15:42"
$ws.Cells.Item($row, 12).Value = 2.75
$ws.Cells.Item($row, 13).Value = "26/10/2023 15:49"
$ws.Cells.Item($row, 14).Value = 3.55
$ws.Cells.Item($row, 15).Value = "22/10/2023 15:42"
$ws.Cells.Item($row, 16).Value = 3.68
$ws.Cells.Item($row, 17).Value = "26/10/2023 15:57"
$ws.Cells.Item($row, 18).Value = 2.15
$ws.Cells.Item($row, 19).Value = "22/10/2023 15:42"
$ws.Cells.Item($row, 20).Value = 2.45
$ws.Cells.Item($row, 21).Value = "26/10/2023 15:57"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/india/isl/north-east-united-jamshedpur/zJzOOWoA/"
